$d = $word.ActiveDocument

# --- Scrum Master name change (paragraph 3 only, so we don't touch the
#     identically-named team member row further down in the table) ---
$pMaster = $d.Paragraphs.Item(3).Range
$pMaster.Find.Execute("Jay Peterson", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Alex Thurston", 2) | Out-Null

# --- Date field: "4/07" -> "03/17" (keeps the underline formatting that
#     was already applied to the digits/slash) ---
$pDate = $d.Paragraphs.Item(4).Range
$pDate.Find.Execute("4/07", $true, $false, $false, $false, $false, `
    $true, 1, $false, "03/17", 2) | Out-Null

# --- Standup table cell content updates ---
$tbl = $d.Tables.Item(1)

$tbl.Cell(2, 2).Range.Find.Execute( `
    "Lots of styling changes for each page.", $true, $false, $false, $false, $false, `
    $true, 1, $false, `
    "Almost completed with wireframe model and required templates, and navigational logic in views.py and urls.py", `
    2) | Out-Null

$tbl.Cell(2, 3).Range.Find.Execute( `
    "Keep updating pages and making them look better. Finalizing UI and resource compatibility.", `
    $true, $false, $false, $false, $false, `
    $true, 1, $false, "Style UI using Bootstrap", 2) | Out-Null

$tbl.Cell(2, 4).Range.Find.Execute( `
    "N/A", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Unfamiliar with Bootstrap.", 2) | Out-Null

$tbl.Cell(3, 2).Range.Find.Execute( `
    "Figured out user authentication and account creation. Figured out event creation.", `
    $true, $false, $false, $false, $false, `
    $true, 1, $false, `
    "Did a little more research but became occupied with other class priorities before vacation.", `
    2) | Out-Null

$tbl.Cell(3, 3).Range.Find.Execute( `
    "Supervisor and lot attendant permissions.", $true, $false, $false, $false, $false, `
    $true, 1, $false, `
    "Vacation " + [char]0x2013 + " will try to complete user authentication procedures after.", `
    2) | Out-Null

$tbl.Cell(3, 4).Range.Find.Execute( `
    "Lots of ways to verify customers. Need to figure out a reliable way.", `
    $true, $false, $false, $false, $false, `
    $true, 1, $false, "Vacation", 2) | Out-Null

$tbl.Cell(4, 2).Range.Find.Execute( `
    "Figured out user profile editing. ", $true, $false, $false, $false, $false, `
    $true, 1, $false, `
    "Refamiliarized self with Django, MVC, data population.", 2) | Out-Null

$tbl.Cell(4, 3).Range.Find.Execute( `
    "Get the password update working. Make html look nice. Unit testing. Verification.", `
    $true, $false, $false, $false, $false, `
    $true, 1, $false, "Begin work on template population", 2) | Out-Null

$tbl.Cell(4, 4).Range.Find.Execute( `
    "Lots of ways to verify customers. Need to figure out a reliable way.", `
    $true, $false, $false, $false, $false, `
    $true, 1, $false, "Other class priorities", 2) | Out-Null

Write-Host "Done."
